$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change
$ws.Range("B1").Value = "VIP Status"

# Row 2 - Barack Obama: swap C2/D2, change G2
$ws.Range("C2").Value = "https://en.wikipedia.org/wiki/Barack_Obama"
$ws.Range("D2").Value = "https://www.instagram.com/barackobama/?hl=en"
$ws.Range("G2").Value = "https://www.obamalibrary.gov/obamas/president-barack-obama"

# Row 3 - Elon Musk: B3 value, shift E/F/G
$ws.Range("B3").Value = 0
$ws.Range("E3").Value = "https://www.forbes.com/profile/elon-musk/"
$ws.Range("F3").Value = "https://www.tesla.com/elon-musk"
$ws.Range("G3").Value = "https://www.spacex.com/"

# Row 4 - Oprah: swap F4/G4 (twitter url case also changes)
$ws.Range("F4").Value = "https://www.facebook.com/oprahwinfrey/"
$ws.Range("G4").Value = "https://twitter.com/oprah"

# Row 5 - Leonardo DiCaprio: E5 add query string
$ws.Range("E5").Value = "https://www.instagram.com/leonardodicaprio/?hl=en"

# Row 6 - Angelina Jolie: E6 instagram post changes
$ws.Range("E6").Value = "https://www.instagram.com/p/CaWTsy5gwQY/?hl=en"

# Row 7 - Jennifer Aniston: swap E7/F7
$ws.Range("E7").Value = "https://www.imdb.com/name/nm0000098/"
$ws.Range("F7").Value = "https://www.instagram.com/p/Czq8K4DMOMQ/?hl=en"

# Row 8 - Robert Downey: D8 drop query string, swap F8/G8
$ws.Range("D8").Value = "https://www.instagram.com/robertdowneyjr/"
$ws.Range("F8").Value = "https://www.facebook.com/robertdowneyjr/"
$ws.Range("G8").Value = "https://twitter.com/robertdowneyjr"

# Row 9 - Scarlett Johansson: G9 changes
$ws.Range("G9").Value = "https://www.nature.com/articles/d41586-024-01578-4"

# Row 11 - paul jarrod frank: B11 value
$ws.Range("B11").Value = 2
